# Apply "update for production and scrap" changes
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "technology" sheet: update production/scrap figures in column B and E
# ---------------------------------------------------------------------------
$tech = $wb.Worksheets.Item("technology")
$tech.Range("B2").Value = 20
$tech.Range("E2").Value = 13
$tech.Range("B3").Value = 20
$tech.Range("E3").Value = 4
$tech.Range("B4").Value = 20
$tech.Range("E4").Value = 13
$tech.Range("B5").Value = 20

# ---------------------------------------------------------------------------
# 2. "renewal" sheet: replace hard-coded / partially formula-driven values
#    with a formula that pulls 30% of the corresponding "capex" sheet cell
# ---------------------------------------------------------------------------
$renewal = $wb.Worksheets.Item("renewal")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($col in $cols) {
    for ($row = 2; $row -le 5; $row++) {
        $cell = "$col$row"
        $renewal.Range($cell).Formula = "=capex!$cell*0.3"
    }
}

# ---------------------------------------------------------------------------
# 3. Update sheet selections (leave the selections the diff records)
# ---------------------------------------------------------------------------
$baseline = $wb.Worksheets.Item("baseline")
$baseline.Range("G2:G14").Select()

$production = $wb.Worksheets.Item("production")
$production.Range("B2:B17").Select()

$renewal.Range("M21").Select()

# Finally activate "technology" and select E5, making it the active tab
$tech.Activate()
$tech.Range("E5").Select()
